# ---------------------------------------------------------------------------
# DU FBS Mock 6 - results update
#
# 1. Switches the "active" sheet from Sheet1 to Sheet3 (tabSelected moves,
#    workbook-level activeTab follows).
# 2. Selects the just-entered data row (A4:XFD4) on Sheet3 so the saved
#    view lands there, matching the "bottomRight" pane selection.
# 3. Fills in the newly-graded answer-key row (row 4, columns C:BZ) on
#    Sheet3 with each student's per-question result, re-using the same
#    "correct" (green) / "wrong" (pink) / "not attempted" (blank) cell
#    styles already used by the sibling rows.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

# Switching the active sheet to Sheet3 is what moves tabSelected="1" off
# Sheet1's sheetView and onto Sheet3's, and bumps workbookView/@activeTab.
$ws3.Activate()

# Reference cells already carrying the three answer-cell styles used
# throughout the sheet (row 5 already has all three), so copying their
# format reuses the existing style indices instead of minting new ones.
$styleGreen = $ws3.Range("C5")    # "correct" answer fill
$stylePink  = $ws3.Range("P5")    # "wrong" answer fill
$styleBlank = $ws3.Range("K5")    # not-attempted / no answer

function SetAnswer($addr, $text, $kind) {
    $ws3.Range($addr).Value = $text
    if ($kind -eq "green") {
        $styleGreen.Copy()
    } else {
        $stylePink.Copy()
    }
    $ws3.Range($addr).PasteSpecial(-4122)
}

function SetBlank($addr) {
    $styleBlank.Copy()
    $ws3.Range($addr).PasteSpecial(-4122)
}

SetAnswer "C4" "C (C)" "green"
SetAnswer "D4" "C (C)" "green"
SetAnswer "E4" "C (C)" "green"
SetAnswer "F4" "C (C)" "green"
SetAnswer "G4" "B (C)" "green"
SetAnswer "H4" "B (C)" "green"
SetAnswer "I4" "A (C)" "green"
SetAnswer "J4" "B (C)" "green"
SetAnswer "K4" "B (C)" "green"
SetAnswer "L4" "C (C)" "green"
SetAnswer "M4" "A (W)" "pink"
SetAnswer "N4" "A (W)" "pink"
SetAnswer "O4" "B (C)" "green"
SetAnswer "P4" "A (W)" "pink"
SetAnswer "Q4" "A (C)" "green"
SetAnswer "R4" "B (C)" "green"
SetAnswer "S4" "C (C)" "green"
SetAnswer "T4" "B (C)" "green"
SetAnswer "U4" "B (C)" "green"
SetAnswer "V4" "A (C)" "green"
SetAnswer "W4" "C (W)" "pink"
SetAnswer "X4" "C (C)" "green"
SetAnswer "Y4" "A (C)" "green"
SetAnswer "Z4" "B (C)" "green"
SetAnswer "AA4" "B (W)" "pink"
SetAnswer "AB4" "D (C)" "green"
SetAnswer "AC4" "B (C)" "green"
SetAnswer "AD4" "C (C)" "green"
SetAnswer "AE4" "C (C)" "green"
SetAnswer "AF4" "A (C)" "green"
SetAnswer "AG4" "C (C)" "green"
SetAnswer "AH4" "B (W)" "pink"
SetAnswer "AI4" "A (W)" "pink"
SetAnswer "AJ4" "D (C)" "green"
SetAnswer "AK4" "A (W)" "pink"
SetAnswer "AL4" "D (C)" "green"
SetBlank  "AM4"
SetBlank  "AN4"
SetAnswer "AO4" "C (W)" "pink"
SetAnswer "AP4" "D (C)" "green"
SetAnswer "AQ4" "B (C)" "green"
SetAnswer "AR4" "A (W)" "pink"
SetAnswer "AS4" "A (C)" "green"
SetAnswer "AT4" "B (W)" "pink"
SetBlank  "AU4"
SetBlank  "AV4"
SetBlank  "AW4"
SetBlank  "AX4"
SetBlank  "AY4"
SetBlank  "AZ4"
SetBlank  "BA4"
SetBlank  "BB4"
SetBlank  "BC4"
SetBlank  "BD4"
SetBlank  "BE4"
SetBlank  "BF4"
SetBlank  "BG4"
SetBlank  "BH4"
SetBlank  "BI4"
SetBlank  "BJ4"
SetAnswer "BK4" "C (C)" "green"
SetAnswer "BL4" "A (C)" "green"
SetAnswer "BM4" "C (C)" "green"
SetAnswer "BN4" "D (W)" "pink"
SetAnswer "BO4" "C (C)" "green"
SetAnswer "BP4" "B (C)" "green"
SetAnswer "BQ4" "B (C)" "green"
SetAnswer "BR4" "C (C)" "green"
SetAnswer "BS4" "B (W)" "pink"
SetAnswer "BT4" "B (C)" "green"
SetAnswer "BU4" "C (C)" "green"
SetAnswer "BV4" "B (C)" "green"
SetAnswer "BW4" "B (C)" "green"
SetBlank  "BX4"
SetAnswer "BY4" "C (W)" "pink"
SetAnswer "BZ4" "C (C)" "green"

# Leave the row we just populated selected, matching the saved view's
# "bottomRight" pane selection (activeCell A4, sqref A4:XFD4).
$ws3.Range("A4:XFD4").Select()

Write-Host "Row 4 answer data written; Sheet3 activated and row selected."
